$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.090.69"
$ws.Range("E2").Value = "  -3.86%  "
$ws.Range("D3").Value = "1.948.40"
$ws.Range("E3").Value = "  -3.72%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'228.05"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -10.45%  "
$ws.Range("D6").Value = "'0.589"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.85%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").Value = "'52.69"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -7.02%  "
$ws.Range("E9").Value = "  -4.67%  "
$ws.Range("D10").Value = "'56.59"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.94%  "
$ws.Range("D11").Value = "'0.0727"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.97%  "
$ws.Range("E12").Value = "  -4.70%  "
$ws.Range("D13").Value = "2.236.37"
$ws.Range("E13").Value = "  -3.70%  "
$ws.Range("D14").Value = "'13.67"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -5.52%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D15").Value = "'0.737"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -9.38%  "
$ws.Range("B16").Value = "Avalanche"
$ws.Range("C16").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D16").Value = "'19.25"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -8.35%  "
$ws.Range("D17").Value = "1.957.54"
$ws.Range("E17").Value = "  -3.65%  "
$ws.Range("E18").Value = "  -7.53%  "
$ws.Range("D19").Value = "36.028.91"
$ws.Range("E19").Value = "  -3.74%  "
$ws.Range("D20").Value = "'66.59"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -4.26%  "
$ws.Range("D21").Value = "0.0₃0782"
$ws.Range("E21").Value = "  -7.78%  "
$ws.Range("D22").Value = "'4.91"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -5.26%  "
$ws.Range("D23").Value = "'218.76"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -4.13%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "'2.29"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -12.32%  "
$ws.Range("D27").Value = "'159.81"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.13%  "
$ws.Range("D28").Value = "'8.34"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -7.43%  "
$ws.Range("D29").Value = "'18.54"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -6.11%  "
$ws.Range("E30").Value = "  -7.78%  "
$ws.Range("D31").Value = "'0.116"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -11.99%  "
$ws.Range("E33").Value = "  -9.24%  "
$ws.Range("D34").Value = "'0.0594"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -11.04%  "
$ws.Range("D35").Value = "'4.14"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -9.14%  "
$ws.Range("B36").Value = "BinanceUSD"
$ws.Range("C36").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.23"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -8.15%  "
$ws.Range("D38").Value = "'1.77"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.50%  "
$ws.Range("D39").Value = "'3.07"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -9.22%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D40").Value = "'2.99"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("B41").Value = "THORChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D41").Value = "'4.95"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -7.35%  "
$ws.Range("D42").Value = "1.391.82"
$ws.Range("E42").Value = "  -1.37%  "
$ws.Range("D43").Value = "'0.0196"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -9.01%  "
$ws.Range("D44").Value = "'0.0850"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -11.71%  "
$ws.Range("E45").Value = "  -12.19%  "
$ws.Range("D46").Value = "'85.41"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -5.66%  "
$ws.Range("D47").Value = "'0.965"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -6.60%  "
$ws.Range("B48").Value = "InjectiveProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D48").Value = "'14.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -9.13%  "
$ws.Range("B49").Value = "MXToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D49").Value = "'2.83"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.33%  "
$ws.Range("E50").Value = "  -8.71%  "
$ws.Range("D51").Value = "2.130.43"
$ws.Range("E51").Value = "  -3.79%  "
